$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.519.97'
$ws.Range('E2').Value = '  -0.32%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.872.56'
$ws.Range('E3').Value = '  -0.46%  '

$ws.Range('E4').Value = '  -2.30%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.25'
$ws.Range('E5').Value = '  -1.23%  '

$ws.Range('E6').Value = '  -1.99%  '

$ws.Range('E7').Value = '  -1.67%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3904'
$ws.Range('E8').Value = '  -1.31%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08370'
$ws.Range('E9').Value = '  +0.27%  '

$ws.Range('E10').Value = '  -1.24%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.80'
$ws.Range('E11').Value = '  -1.04%  '

$ws.Range('E12').Value = '  -1.01%  '

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.875.41'
$ws.Range('E13').Value = '  +1.21%  '

$ws.Range('B14').Value = 'Solana'
$ws.Range('C14').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.48'
$ws.Range('E14').Value = '  -0.62%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.275'
$ws.Range('E15').Value = '  +0.18%  '

$ws.Range('E16').Value = '  -2.35%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001101'
$ws.Range('E17').Value = '  -1.06%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '91.23'
$ws.Range('E18').Value = '  -0.38%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06733'
$ws.Range('E19').Value = '  -0.94%  '

$ws.Range('E20').Value = '  -0.12%  '

$ws.Range('E21').Value = '  -2.01%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.926'
$ws.Range('E22').Value = '  -1.29%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.548.97'
$ws.Range('E23').Value = '  -0.36%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.10'
$ws.Range('E24').Value = '  -0.86%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.189'
$ws.Range('E25').Value = '  -4.01%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.080.04'
$ws.Range('E26').Value = '  +0.73%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '158.67'
$ws.Range('E27').Value = '  -2.37%  '

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.63'
$ws.Range('E28').Value = '  -1.14%  '

$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.433'
$ws.Range('E29').Value = '  +1.98%  '

$ws.Range('E30').Value = '  -0.48%  '

$ws.Range('E31').Value = '  -1.53%  '

$ws.Range('E32').Value = '  +0.78%  '

$ws.Range('E33').Value = '  -1.96%  '

$ws.Range('E34').Value = '  -1.33%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02459'
$ws.Range('E35').Value = '  +0.77%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06611'
$ws.Range('E36').Value = '  +1.16%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2169'
$ws.Range('E37').Value = '  -0.82%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.915'
$ws.Range('E38').Value = '  -3.27%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.039'
$ws.Range('E39').Value = '  +0.49%  '

$ws.Range('E40').Value = '  -0.68%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.239'
$ws.Range('E41').Value = '  -1.12%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6388'
$ws.Range('E42').Value = '  -1.31%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.11'
$ws.Range('E43').Value = '  -0.97%  '

$ws.Range('E44').Value = '  -1.88%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6015'
$ws.Range('E45').Value = '  -0.80%  '

$ws.Range('E46').Value = '  +0.35%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.679'
$ws.Range('E47').Value = '  -1.21%  '

$ws.Range('E48').Value = '  +0.44%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.214'
$ws.Range('E49').Value = '  -0.16%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '122.64'
$ws.Range('E50').Value = '  +0.29%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06806'
$ws.Range('E51').Value = '  -1.00%  '
